$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the close price on the existing last row (306): 30.75 -> 30.85
$ws.Range("F306").Value = 30.85

# Append three new data rows (307-309), copying row 306's formatting down first
# so the new rows inherit the same styles (e.g. the date style on column A).
$ws.Range("A306:G306").Copy($ws.Range("A307:G307"))
$ws.Range("A306:G306").Copy($ws.Range("A308:G308"))
$ws.Range("A306:G306").Copy($ws.Range("A309:G309"))

# Row 307
$ws.Range("A307").Value = 45047.33333333334
$ws.Range("B307").Value = "FX_IDC:USDEGP"
$ws.Range("C307").Value = 30.9499
$ws.Range("D307").Value = 30.9499
$ws.Range("E307").Value = 30.73
$ws.Range("F307").Value = 30.85
$ws.Range("G307").Value = 0

# Row 308
$ws.Range("A308").Value = 45078.33333333334
$ws.Range("B308").Value = "FX_IDC:USDEGP"
$ws.Range("C308").Value = 30.85
$ws.Range("D308").Value = 30.9499
$ws.Range("E308").Value = 30.75
$ws.Range("F308").Value = 30.85
$ws.Range("G308").Value = 0

# Row 309
$ws.Range("A309").Value = 45110.33333333334
$ws.Range("B309").Value = "FX_IDC:USDEGP"
$ws.Range("C309").Value = 30.85
$ws.Range("D309").Value = 30.9499
$ws.Range("E309").Value = 30.75
$ws.Range("F309").Value = 30.83
$ws.Range("G309").Value = 0
